# "Tried to implement Penalty Reward System (unfinished)"
# Shifts the Forecast Comparison week-start dates forward by one week and
# updates the MyForecast values, then refreshes the dependent Summary stats.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Keep the date-looking / numeric-looking strings stored as TEXT (matching the
# original inline-string cells) instead of letting Excel auto-convert them to
# real dates/numbers.
$wsForecast.Range("B2:B17").NumberFormat = "@"
$wsSummary.Range("B2:B15").NumberFormat = "@"

# --- Forecast Comparison sheet: Week_Start_Date (B) + MyForecast (D) ---

$wsForecast.Range("B2").Value = "2025-01-12"
$wsForecast.Range("D2").Value = 62

$wsForecast.Range("B3").Value = "2025-01-19"
$wsForecast.Range("D3").Value = 68

$wsForecast.Range("B4").Value = "2025-01-26"
$wsForecast.Range("D4").Value = 67

$wsForecast.Range("B5").Value = "2025-02-02"
$wsForecast.Range("D5").Value = 68

$wsForecast.Range("B6").Value = "2025-02-09"
$wsForecast.Range("D6").Value = 69

$wsForecast.Range("B7").Value = "2025-02-16"
$wsForecast.Range("D7").Value = 70

$wsForecast.Range("B8").Value = "2025-02-23"
$wsForecast.Range("D8").Value = 71

$wsForecast.Range("B9").Value = "2025-03-02"
$wsForecast.Range("D9").Value = 70

$wsForecast.Range("B10").Value = "2025-03-09"
$wsForecast.Range("D10").Value = 69

$wsForecast.Range("B11").Value = "2025-03-16"
$wsForecast.Range("D11").Value = 68

$wsForecast.Range("B12").Value = "2025-03-23"
$wsForecast.Range("D12").Value = 67

$wsForecast.Range("B13").Value = "2025-03-30"
$wsForecast.Range("D13").Value = 67

$wsForecast.Range("B14").Value = "2025-04-06"
$wsForecast.Range("D14").Value = 68

$wsForecast.Range("B15").Value = "2025-04-13"
$wsForecast.Range("D15").Value = 69

$wsForecast.Range("B16").Value = "2025-04-20"
$wsForecast.Range("D16").Value = 69

$wsForecast.Range("B17").Value = "2025-04-27"
$wsForecast.Range("D17").Value = 69

# --- Summary sheet: recomputed stats from the new MyForecast column ---

$wsSummary.Range("B2").Value  = "2023-01-22 to 2025-01-05"
$wsSummary.Range("B4").Value  = "199"
$wsSummary.Range("B5").Value  = "54"
$wsSummary.Range("B6").Value  = "43"
$wsSummary.Range("B7").Value  = "47"
$wsSummary.Range("B8").Value  = "4266 units"
$wsSummary.Range("B9").Value  = "1092"
$wsSummary.Range("B10").Value = "545"
$wsSummary.Range("B11").Value = "264"
$wsSummary.Range("B12").Value = "71"
$wsSummary.Range("B14").Value = "62"
$wsSummary.Range("B15").Value = "2025-01-12"
